# Adds Fall 2022 Week 6 matchup data rows (1242-1312) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(4,0,4,3),
    @(2,2,2,1),
    @(5,0,3,2),
    @(6,0,7,2),
    @(6,0,5,2),
    @(5,0,6,2),
    @(5,0,4,2),
    @(5,0,5,2),
    @(4,3,3,0),
    @(4,2,4,1),
    @(4,0,5,2),
    @(4,3,3,0),
    @(5,2,5,1),
    @(6,2,6,0),
    @(2,2,2,1),
    @(5,2,5,0),
    @(7,2,6,1),
    @(5,2,5,0),
    @(3,0,4,3),
    @(2,2,4,0),
    @(6,2,7,0),
    @(2,0,3,3),
    @(6,1,4,2),
    @(3,3,2,0),
    @(7,2,5,0),
    @(3,2,4,0),
    @(5,1,5,2),
    @(3,2,6,0),
    @(4,0,4,3),
    @(6,2,6,1),
    @(5,2,4,0),
    @(4,0,5,3),
    @(4,3,4,0),
    @(3,0,2,3),
    @(3,0,2,2),
    @(5,3,5,0),
    @(4,0,3,3),
    @(6,1,7,2),
    @(5,0,6,3),
    @(2,2,2,1),
    @(5,1,6,2),
    @(5,1,5,2),
    @(4,0,2,2),
    @(6,2,6,0),
    @(2,3,2,0),
    @(5,0,6,2),
    @(4,2,3,1),
    @(6,2,6,0),
    @(4,1,3,2),
    @(4,1,5,2),
    @(4,2,4,0),
    @(6,0,6,3),
    @(5,3,3,0),
    @(2,3,3,0),
    @(3,3,4,0),
    @(3,2,4,1),
    @(5,2,5,0),
    @(7,2,6,1),
    @(3,3,3,0),
    @(5,2,5,1),
    @(3,2,3,1),
    @(6,2,4,0),
    @(6,2,6,0),
    @(3,1,4,2),
    @(5,1,5,2),
    @(6,0,6,2),
    @(6,0,4,2),
    @(4,2,7,0),
    @(6,2,6,0),
    @(5,3,4,0),
    @(4,2,4,1)
)

$startRow = 1242
$endRow = $startRow + $newData.Count - 1

$arr = New-Object 'object[,]' $newData.Count,4
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i, $j] = $rowVals[$j]
    }
}

$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$targetRange.Value = $arr

$nextCell = $ws.Cells.Item($endRow + 1, 1)
$nextCell.Select()
